$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the "Inf" index column (A): it was off by one (1..5 instead of 0..4) ---
# Shift A2:A6 values down into A3:A7 (i.e. each label moves one row down,
# keeping the same numeric value), and clear the now-unused A2 cell.
# Process bottom-up so we don't clobber a value before reading it.
$ws.Range("A7").Value = $ws.Range("A6").Value()
$ws.Range("A6").Value = $ws.Range("A5").Value()
$ws.Range("A5").Value = $ws.Range("A4").Value()
$ws.Range("A4").Value = $ws.Range("A3").Value()
$ws.Range("A3").Value = $ws.Range("A2").Value()
$ws.Range("A2").ClearContents()

# --- Fix row 2: the sensitivity values had drifted one column to the right ---
# Each of C2:F2 should hold what used to be one column to its left (B2:E2),
# and the old F2 value is discarded since it no longer belongs to the series.
# Process right-to-left so we don't clobber a value before reading it.
$ws.Range("F2").Value = $ws.Range("E2").Value()
$ws.Range("E2").Value = $ws.Range("D2").Value()
$ws.Range("D2").Value = $ws.Range("C2").Value()
$ws.Range("C2").Value = $ws.Range("B2").Value()

# B2 is now derived from the corrected C2 (C2 - 2), matching the rest of the series
$ws.Range("B2").Formula = "=C2-2"
